# US2: Casos de prueba completados
#
# 1. Extend the "expected value" text of step 3 in the "Pasos" sheet.
# 2. Make "Pasos" the active/selected sheet (was "Control de cambios").
# 3. Update the cell selection on both sheets accordingly.

$wb = $excel.ActiveWorkbook

$pasos = $wb.Worksheets.Item("Pasos")
$cambios = $wb.Worksheets.Item("Control de cambios")

# Update the expected-result text for step 3 (row 4, column C) on "Pasos".
$pasos.Range("C4").Value = "Se carga la pagina BuscarPlayas, con todas las playas de <Ciudad1> disponibles en un mapa. Se carga la informacion de las playas en la grilla de playas debajo del mapa."

# Activate the "Pasos" sheet (making it the selected/active tab) and move
# the selection there to C5. "Control de cambios" keeps its previous
# selection (B3), so it is left untouched.
$pasos.Activate()
$pasos.Range("C5").Select()
